$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '$2b$12$5S12DDSfSY1d7mYVvJ8CBetURqyTzPtb3QBjH8hHSMwRRI.MnjP/2'
$ws.Range("B3").Value = '$2b$12$W/Ro6oHLXLgtZ2cL.xMslOjXkt2hPRGSiym88Kli9jtqJeLrrwsh6'
$ws.Range("B4").Value = '$2b$12$9tmFMt9l/HZgYIZ5DWVM8.6CSSEU4zRgNBEIoe/o/H0oxY4pP5bcK'
$ws.Range("B5").Value = '$2b$12$x/1LHsLV2xKTfn915sPtb..DZIGsR9eEx7fftzyipz9XhH1E093oe'
